$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force the new cells to be treated as text so "1000.0" is preserved
# literally instead of being normalized into the number 1000.
$ws.Range("A2:B2").NumberFormat = "@"

$ws.Range("A2").Value = "1000.0"
$ws.Range("B2").Value = "1000.0"

# Drop the explicit number-format styling again so the new cells end up
# without any style index, matching the original sheet's plain cells.
$ws.Range("A2:B2").ClearFormats()
